# Updated files from WRI China bringing the Hong Kong EPS up to v2.0.0
#
# Changes applied to "Elasticity of Dist Solar Deployment wrt Subsidy Perc.xlsx":
#   - About sheet: add a "Notes:" section explaining the elasticities
#   - EoDSDwSP sheet: clarify the "Elasticity" label as "Elasticity (dimensionless)"

$wb = $excel.ActiveWorkbook

# --- "About" sheet: append explanatory notes below the existing source info ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A9").Value  = "Notes:"
$wsAbout.Range("A10").Value = "Elasticities intended to reflect change in deployment with changing"
$wsAbout.Range("A11").Value = "distributed solar price (through subsidies)."

# --- "EoDSDwSP" sheet: clarify the elasticity header, move selection ---
$wsResult = $wb.Worksheets.Item("EoDSDwSP")
$wsResult.Range("B1").Value = "Elasticity (dimensionless)"
$wsResult.Range("B2").Select()

# Leave the "About" sheet as the active tab/selection, as in the source workbook
$wsAbout.Activate()
$wsAbout.Range("A12").Select()
